$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1323.8889
$ws.Range("J70").Value = 1239.375
$ws.Range("L70").Value = 3718.125
$ws.Range("N70").Value = -4258.125
$ws.Range("H73").Value = 1323.8889
$ws.Range("J73").Value = 1239.375
$ws.Range("L73").Value = 3718.125
$ws.Range("N73").Value = -5590.125
$ws.Range("H112").Value = 2930.3845
$ws.Range("I112").Value = 1398
$ws.Range("J112").Value = 3058.0833
$ws.Range("K112").Value = 4194
$ws.Range("L112").Value = 9174.249899999999
$ws.Range("M112").Value = -3086
$ws.Range("N112").Value = -11390.2499
$ws.Range("H121").Value = 856
$ws.Range("J121").Value = 856
$ws.Range("L121").Value = 2568
$ws.Range("N121").Value = -6062
$ws.Range("H132").Value = 2867.739
$ws.Range("I132").Value = 1559.875
$ws.Range("K132").Value = 4679.625
$ws.Range("M132").Value = -2149.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 193.5
$ws.Range("I3").Value = 193.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 193.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -78.5
$ws.Range("N3").ClearContents()
$ws.Range("H32").Value = 4381.26
$ws.Range("I32").Value = 3313.9167
$ws.Range("K32").Value = 3313.9167
$ws.Range("M32").Value = -3026.9167
$ws.Range("H122").Value = 1775.375
$ws.Range("I122").Value = 1775.375
$ws.Range("K122").Value = 5326.125
$ws.Range("M122").Value = -2876.125
$ws.Range("H132").Value = 2700.3125
$ws.Range("I132").Value = 1921.3334
$ws.Range("K132").Value = 5764.0002
$ws.Range("M132").Value = -3234.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H134").Value = 2575.75
$ws.Range("I134").Value = 2356.6
$ws.Range("K134").Value = 7069.799999999999
$ws.Range("M134").Value = -4534.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 79976.7
$ws.Range("J18").Value = 79976.7
$ws.Range("L18").Value = 79976.7
$ws.Range("N18").Value = -80436.7
$ws.Range("H28").Value = 14748.5
$ws.Range("J28").Value = 14748.5
$ws.Range("L28").Value = 14748.5
$ws.Range("N28").Value = -15238.5
$ws.Range("H43").Value = 39666.332
$ws.Range("J43").Value = 39666.332
$ws.Range("L43").Value = 39666.332
$ws.Range("N43").Value = -40034.332
$ws.Range("H56").Value = 49999
$ws.Range("J56").Value = 49999
$ws.Range("L56").Value = 49999
$ws.Range("N56").Value = -51689
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20406
$ws.Range("H88").Value = 15797.223
$ws.Range("J88").Value = 15797.223
$ws.Range("L88").Value = 15797.223
$ws.Range("N88").Value = -16609.223
$ws.Range("H91").Value = 15797.223
$ws.Range("J91").Value = 15797.223
$ws.Range("L91").Value = 15797.223
$ws.Range("N91").Value = -18605.223
$ws.Range("H101").Value = 39666.332
$ws.Range("J101").Value = 39666.332
$ws.Range("L101").Value = 39666.332
$ws.Range("N101").Value = -46156.332
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H134").Value = 3413.182
$ws.Range("I134").Value = 3394
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 10182
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -7647
$ws.Range("N134").Value = -15568.5
$ws.Range("H136").Value = 20000
$ws.Range("J136").Value = 20000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -65100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 345.27274
$ws.Range("I8").Value = 345.27274
$ws.Range("K8").Value = 1035.81822
$ws.Range("M8").Value = -896.8182200000001
$ws.Range("H46").Value = 98.5
$ws.Range("I46").Value = 98.5
$ws.Range("K46").Value = 295.5
$ws.Range("M46").Value = -204.5
$ws.Range("H109").Value = 723.44446
$ws.Range("I109").Value = 623
$ws.Range("J109").Value = 849
$ws.Range("K109").Value = 1869
$ws.Range("L109").Value = 2547
$ws.Range("M109").Value = -829
$ws.Range("N109").Value = -4627
$ws.Range("H128").Value = 499994
$ws.Range("I128").Value = 499994
$ws.Range("K128").Value = 1499982
$ws.Range("M128").Value = -1495002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16679299
$ws.Range("I70").Value = 16679299
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 16679299
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -16679029
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 16679299
$ws.Range("I73").Value = 16679299
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 16679299
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -16678363
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 3006.7778
$ws.Range("J80").Value = 3306.4
$ws.Range("L80").Value = 3306.4
$ws.Range("N80").Value = -5302.4
$ws.Range("H83").Value = 3006.7778
$ws.Range("J83").Value = 3306.4
$ws.Range("L83").Value = 16532
$ws.Range("N83").Value = -26516
$ws.Range("H99").Value = 9030
$ws.Range("I99").Value = 9030
$ws.Range("K99").Value = 9030
$ws.Range("M99").Value = -6784
$ws.Range("H102").Value = 663.2
$ws.Range("J102").Value = 1000
$ws.Range("L102").Value = 1000
$ws.Range("N102").Value = -4244
$ws.Range("H126").Value = 1000000000
$ws.Range("I126").Value = 1000000000
$ws.Range("K126").Value = 3000000000
$ws.Range("M126").Value = -2999997530
$ws.Range("H132").Value = 2855.652
$ws.Range("I132").Value = 2430.625
$ws.Range("J132").Value = 3827.1428
$ws.Range("K132").Value = 7291.875
$ws.Range("L132").Value = 11481.4284
$ws.Range("M132").Value = -4761.875
$ws.Range("N132").Value = -16541.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 30333.334
$ws.Range("I25").Value = 30000
$ws.Range("K25").Value = 30000
$ws.Range("M25").Value = -29770
$ws.Range("H40").Value = 11933
$ws.Range("I40").Value = 5899.5
$ws.Range("J40").Value = 24000
$ws.Range("K40").Value = 5899.5
$ws.Range("L40").Value = 24000
$ws.Range("M40").Value = -5763.5
$ws.Range("N40").Value = -24272
$ws.Range("H46").Value = 3014.4783
$ws.Range("I46").Value = 1911
$ws.Range("J46").Value = 3723.8572
$ws.Range("K46").Value = 1911
$ws.Range("L46").Value = 3723.8572
$ws.Range("M46").Value = -1723
$ws.Range("N46").Value = -4099.8572
$ws.Range("H68").Value = 2600.5
$ws.Range("I68").Value = 2198
$ws.Range("J68").Value = 3003
$ws.Range("K68").Value = 2198
$ws.Range("L68").Value = 3003
$ws.Range("M68").Value = -1449
$ws.Range("N68").Value = -4501
$ws.Range("H71").Value = 2600.5
$ws.Range("I71").Value = 2198
$ws.Range("J71").Value = 3003
$ws.Range("K71").Value = 10990
$ws.Range("L71").Value = 15015
$ws.Range("M71").Value = -7246
$ws.Range("N71").Value = -22503
$ws.Range("H136").Value = 5342.7144
$ws.Range("I136").Value = 5099.8
$ws.Range("K136").Value = 15299.4
$ws.Range("M136").Value = -12749.4
